# Generate Report for handoff
# - The old source file "9b7c72f4-0680-49fc-82db-c9ff185fce14.md" was renamed/replaced by
#   "94881592-d16f-4818-a514-0e72c6c88ad8.md" and got a fresh handoff (new hash, new datetime).
# - The row describing "8bf6ccdf-ad48-4cf9-a358-992f610c9b85.md" / "Handoff failed" is removed entirely
#   from every sheet (the file no longer participates), so what used to be row 4 (.localization-config)
#   becomes row 3.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------------
# Sheet 1: "Overview"
# ----------------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Drop all existing hyperlinks up front (will be rebuilt below with the final content/targets)
$ws1.Hyperlinks.Delete()

# Remove the whole "8bf6ccdf..." row (old row 3); row 4 (.localization-config) shifts up to row 3
$ws1.Rows.Item(3).Delete()

# Row 2: source file renamed
$ws1.Range("A2").Value = "94881592-d16f-4818-a514-0e72c6c88ad8.md"
$ws1.Range("B2").Value = "Not yet handed off"
$ws1.Range("C2").Value = "Not yet handed off"

# Row 3 (previously row 4): .localization-config, unaffected content
$ws1.Range("A3").Value = ".localization-config"
$ws1.Range("B3").Value = "Not localized"
$ws1.Range("C3").Value = "Not localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f89a78dcd934b12371db4bd5f34af7910f7a52ab/e2e/94881592-d16f-4818-a514-0e72c6c88ad8.md", [Type]::Missing, [Type]::Missing, "94881592-d16f-4818-a514-0e72c6c88ad8.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f89a78dcd934b12371db4bd5f34af7910f7a52ab/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ----------------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ----------------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

# Remove the whole "8bf6ccdf..." row (old row 3); row 4 (.localization-config) shifts up to row 3
$ws2.Rows.Item(3).Delete()

# Row 2: renamed source file, fresh handoff file/datetime
$ws2.Range("A2").Value = "94881592-d16f-4818-a514-0e72c6c88ad8.md"
$ws2.Range("B2").Value = "Not yet handed off"
$ws2.Range("C2").Value = "94881592-d16f-4818-a514-0e72c6c88ad8.add684ab223b83648bfc9fd15f9735620d7a49ef.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-08 14:49:12"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"

# Row 3 (previously row 4): .localization-config, unaffected content
$ws2.Range("A3").Value = ".localization-config"
$ws2.Range("B3").Value = "Not localized"
$ws2.Range("D3").Value = "0001-01-01 00:00:00"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f89a78dcd934b12371db4bd5f34af7910f7a52ab/e2e/94881592-d16f-4818-a514-0e72c6c88ad8.md", [Type]::Missing, [Type]::Missing, "94881592-d16f-4818-a514-0e72c6c88ad8.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/76aeb29eb83b51894690e59da6d998f7f81bc1c8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/94881592-d16f-4818-a514-0e72c6c88ad8.add684ab223b83648bfc9fd15f9735620d7a49ef.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "94881592-d16f-4818-a514-0e72c6c88ad8.add684ab223b83648bfc9fd15f9735620d7a49ef.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f89a78dcd934b12371db4bd5f34af7910f7a52ab/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ----------------------------------------------------------------------------------
# Sheet 3: "de-de"
# ----------------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

# Remove the whole "8bf6ccdf..." row (old row 3); row 4 (.localization-config) shifts up to row 3
$ws3.Rows.Item(3).Delete()

# Row 2: renamed source file, fresh handoff file/datetime
$ws3.Range("A2").Value = "94881592-d16f-4818-a514-0e72c6c88ad8.md"
$ws3.Range("B2").Value = "Not yet handed off"
$ws3.Range("C2").Value = "94881592-d16f-4818-a514-0e72c6c88ad8.add684ab223b83648bfc9fd15f9735620d7a49ef.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-08 14:49:26"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"

# Row 3 (previously row 4): .localization-config, unaffected content
$ws3.Range("A3").Value = ".localization-config"
$ws3.Range("B3").Value = "Not localized"
$ws3.Range("D3").Value = "0001-01-01 00:00:00"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f89a78dcd934b12371db4bd5f34af7910f7a52ab/e2e/94881592-d16f-4818-a514-0e72c6c88ad8.md", [Type]::Missing, [Type]::Missing, "94881592-d16f-4818-a514-0e72c6c88ad8.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/11574741c907abb284053cac43d48a4170f32011/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/94881592-d16f-4818-a514-0e72c6c88ad8.add684ab223b83648bfc9fd15f9735620d7a49ef.de-de.xlf", [Type]::Missing, [Type]::Missing, "94881592-d16f-4818-a514-0e72c6c88ad8.add684ab223b83648bfc9fd15f9735620d7a49ef.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f89a78dcd934b12371db4bd5f34af7910f7a52ab/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
